$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels: underscores -> spaces (D1, G1, H1); I1 keeps its value.
$ws.Range("D1").Value = "sales price"
$ws.Range("G1").Value = "bullet points"
$ws.Range("H1").Value = "contact us"

# Move the active selection to I1 (also clears the old topLeftCell/selection scroll state).
[void]$ws.Range("I1").Select()
